$wb = $excel.ActiveWorkbook

# --- "Clasificacion productos": append two rows to the bottom of the
# Linea de producto / Producto list (Almacen / yogures_naturales_y_sabores
# and Almacen / yogures_y_postres_infantiles) ---
$wsProd = $wb.Worksheets.Item("Clasificacion productos")
$wsProd.Activate()
$wsProd.Range("B150").Value = "Almacen"
$wsProd.Range("C150").Value = "yogures_naturales_y_sabores"
$wsProd.Range("B151").Value = "Almacen"
$wsProd.Range("C151").Value = "yogures_y_postres_infantiles"
$wsProd.Range("C154").Select()

# --- "medios de pago": append two rows (Ewallet/Billetera Electronica and
# Cash/Efectivo) repeating the existing payment methods ---
$wsPago = $wb.Worksheets.Item("medios de pago")
$wsPago.Activate()
$wsPago.Range("B6").Value = "Ewallet"
$wsPago.Range("C6").Value = "Billetera Electronica"
$wsPago.Range("B7").Value = "Cash"
$wsPago.Range("C7").Value = "Efectivo"
$wsPago.Range("B7:C7").Select()
